$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 21, 22, 23, 24, 25, 26)
$bValues = @(
    'eficazmente, unvaccinated, incuráveis, responsabilizadas, tetyana, barré, gardasil, danificam, isentando, obukhanych, triplice, bulas, esboçam, patrocinados, –, prevalente, calá, lesadas, conservantes, razoes, comprovadamente, monossódico, vaxtruth, vaers, sgb, ​​e, nauseam, —, sanevax, guillain, empurradores, imuno, vactruth, afligidos, vacinadas, ​​em, vactruthcom, debilitados, matá, anvisa, adjuvantes, glutamato, concedem, ’, concertado, empurrá, saudáveis​​, contraia, imunocomprometidas, “, v…/…, gms, placebos, sids, previnem, imunológica, reorientar',
    'linfomas, acidificando, entopem, dosagens, aftas, albicans, piccoli, catedráticos, mês1, metastasipolmonari, peritonial, novembro/2008, aproveitável, administrá, jfj, complementá, candidíase, peritoneale, colecisti, conjuntiva, carcinosi, metástases, 360°, coliciste, cancerologia, nahco3, epatocarcinoma, dellintestino, adaptabilidade, trattati, prostata, midollare, pecezinho, mangueirinha, 90°, gotejador, hepatocarcinoma, axilares, revolucionando, aprofundando, palpebra, cancerosas, penetrarem, embasamento, tumore, cérvico, muitíssimos, polmonari, bórax, cerebrali, semanas4, deleto, removê, charlatanismo, diretaço, câncer/2008, vescica, endossando, hodgkin, didática, ½, uterino, persistirem, salgar, alcuni, difusível, ingerível, casi, empiricamente, metodologias, especificas, subsequentemente, metodicamente, coróide, urinaria, marcosrochablogspotcom//descoberta, diffuso, adenocarcinoma, vezes1, oriundos, inalações, reversões, teclar, 10/12/08, arteriografia, percebíamos, candidíases, apregoam, simoncini, bronchiale, 5todos, brônquios, equimoses, integradores, alcalinizando, ministrando, homeopático, escorra, intracelular, •, 1/2, turmores, praticadas, cobaltoterapia, menstruações, linfonodos, tanti, subministrar, intertítulo, intersticial, polmoni, paradigmas, picaretagem, palpáveis, 36a, orientadas, poderosíssima, /100, cervice, obtinha, constatou, homeopatia, lincados, dias3, tumori, “, fungophp, oncologistas, macrobióticos, vari',
    'mg/929, “, incrementando, aspirará, cm2, 16º, especificam, aceitável…, repassem',
    'fibromialgia, aspartame, agravando, comercializem, adiposas, assintomáticas, internadas, colocá, riccio, 30º, engaje, mancy, fórmico, fetais, nutrasweet, dietética, assolassem, arckle, malformações, dietético, suplanta, sistêmico, articulares, entitulado, zerocal',
    'partilhem',
    'aquecê, progrida, estreitarem, lugarquando, esvaziamos',
    'menstruais, solidifica, obstipação, vómitos, reparte',
    '–, tamiflu, nossso, acerola, h1n1, famíliares, 12/horas, infectologista',
    'cardiologia/secretaria, socbrasde, fenofinol, teixeira, ambev, indenizar, “, freitas, ingerirem, usp, almeido, repassado, conseqüente, voliteral',
    'proprios, socbrasde, cardiologia/secretaria, voliteral, fenofinol, teixeira, indenizar, repassada, “, freitas, ingerirem, substancias, usp, almeido, conseqüente, repassem',
    'osh, encaminhá, cancerosa, cancerosas, curativas, incentivei, comprovaram, digeri, handphone, polifenol, guruprasad',
    'alforjón, b17, laetril, sabugueiro, sorgo, documente, macadâmia, “, grumos, rebentos, amígdalina, absorvível, 5sementes, tonsilar, milheto, colonizador',
    '11h, indicativos, 16h',
    'lpki',
    'subornam, distribuíram, partilhe',
    'superdivertida, cevabacillus, muniz, fatie, ativus, –, conservantes, internamento, divulguem, camboriú, coccus, estabilizantes, sacarovictus, ‘, 03/01/16, 08/01/16, ’, fatiados, minimiza, dermato',
    '80cc, coronavírus, 50cc, invadirá',
    'medicos, afectados, informaçao, adiquire, pergosa, japao, transmissao',
    'globalista, totall, começõu, wurhan, elevem, jinping, aglomerações, rockfellers, globalistas, soros, gloso, rothschilds, 2/3',
    'economizaremos, arrastrar, circulante, executá, escarro, 19h25, enviá, ​​e',
    'compartilhá, reencaminhado',
    'partilhou, ademola',
    'coronavírus, wuhan, aconselhavam, repatriados',
    'alcalina, anticâncer, huiren, polifenóis, “, enfatizou, cistos'
)
$cValues = @(
    0.03988803358992302,
    0.03083491461100569,
    0.02593659942363112,
    0.03019323671497584,
    0.0101010101010101,
    0.01633986928104575,
    0.008183306055646482,
    0.05128205128205128,
    0.05857740585774059,
    0.05614035087719298,
    0.01906412478336222,
    0.02531645569620253,
    0.015,
    0.003759398496240601,
    0.03333333333333333,
    0.04608294930875576,
    0.0196078431372549,
    0.0958904109589041,
    0.0755813953488372,
    0.0199501246882793,
    0.006042296072507553,
    0.01162790697674419,
    0.01801801801801802,
    0.03043478260869565
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
    $ws.Cells.Item($r, 3).Value = $cValues[$i]
}
